$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111396053
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "80"
$ws.Range("P2").Value = "S om järnvägen - 2, Vg"
$ws.Range("Q2").Value = 432083.280685614
$ws.Range("R2").Value = 6419676.539718015

# Row 3
$ws.Range("A3").Value = 111396060
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "90"
$ws.Range("J3").Value = ""
$ws.Range("P3").Value = "S om järnvägen - 3, Vg"
$ws.Range("Q3").Value = 432076.641898193
$ws.Range("R3").Value = 6419661.774153749

# Row 4
$ws.Range("A4").Value = 111396045
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "70"
$ws.Range("J4").Value = "stjälkar/strån/skott"
$ws.Range("P4").Value = "S om järnvägen, Vg"
$ws.Range("Q4").Value = 431889.3909100805
$ws.Range("R4").Value = 6419670.266848063

# Row 6
$ws.Range("A6").Value = 111482955
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "70"
$ws.Range("P6").Value = "S om järnvägen - 5, Vg"
$ws.Range("Q6").Value = 432064.1298546481
$ws.Range("R6").Value = 6419677.395781181

# Row 7
$ws.Range("A7").Value = 111482980
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "10"
$ws.Range("P7").Value = "S om järnvägen - 6, Vg"
$ws.Range("Q7").Value = 432048.2263952638
$ws.Range("R7").Value = 6419681.385014677

# Row 8
$ws.Range("A8").Value = 111483105
$ws.Range("P8").Value = "S om järnvägen - 8, Vg"
$ws.Range("Q8").Value = 431947.1499479365
$ws.Range("R8").Value = 6419623.056550305

# Row 10
$ws.Range("A10").Value = 111491187
$ws.Range("P10").Value = "S om järnvägen - 18, Vg"
$ws.Range("Q10").Value = 431829.514510141
$ws.Range("R10").Value = 6419749.394753682

# Row 11
$ws.Range("A11").Value = 111490843
$ws.Range("B11").Value = 96348
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = "Knärot"
$ws.Range("G11").Value = "Goodyera repens"
$ws.Range("H11").Value = "(L.) R. Br."
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "50"
$ws.Range("J11").Value = "stjälkar/strån/skott"
$ws.Range("K11").Value = "fullt utvecklade blad"
$ws.Range("L11").Value = ""
$ws.Range("P11").Value = "S om järnvägen - 17, Vg"
$ws.Range("Q11").Value = 431803.2980747336
$ws.Range("R11").Value = 6419679.170503675
$ws.Range("AJ11").ClearContents() | Out-Null
$ws.Range("AK11").ClearContents() | Out-Null
$ws.Range("AM11").ClearContents() | Out-Null
$ws.Range("AO11").ClearContents() | Out-Null

# Row 12
$ws.Range("A12").Value = 111483197
$ws.Range("B12").Value = 73689
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 308
$ws.Range("F12").Value = "Brunpudrad nållav"
$ws.Range("G12").Value = "Chaenotheca gracillima"
$ws.Range("H12").Value = "(Vain.) Tibell"
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").ClearContents() | Out-Null
$ws.Range("P12").Value = "S om järnvägen - 11, Vg"
$ws.Range("Q12").Value = 431937.082796899
$ws.Range("R12").Value = 6419625.884406033
$ws.Range("AM12").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO12").Value = "Standing dead tree/snags"

# Row 13
$ws.Range("A13").Value = 111483037
$ws.Range("B13").Value = 96348
$ws.Range("D13").Value = "VU"
$ws.Range("E13").Value = 220787
$ws.Range("F13").Value = "Knärot"
$ws.Range("G13").Value = "Goodyera repens"
$ws.Range("H13").Value = "(L.) R. Br."
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "60"
$ws.Range("J13").Value = "stjälkar/strån/skott"
$ws.Range("K13").Value = "blomning"
$ws.Range("L13").Value = ""
$ws.Range("P13").Value = "S om järnvägen - 7, Vg"
$ws.Range("Q13").Value = 432060.6482816387
$ws.Range("R13").Value = 6419660.45125766
$ws.Range("AJ13").ClearContents() | Out-Null
$ws.Range("AK13").ClearContents() | Out-Null
$ws.Range("AM13").ClearContents() | Out-Null
$ws.Range("AO13").ClearContents() | Out-Null

# Row 14
$ws.Range("A14").Value = 111483437
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "100"
$ws.Range("P14").Value = "S om järnvägen - 15, Vg"
$ws.Range("Q14").Value = 431797.479853621
$ws.Range("R14").Value = 6419681.394993878

# Row 16
$ws.Range("A16").Value = 111483107
$ws.Range("B16").Value = 73681
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 6439
$ws.Range("F16").Value = "Gulnål"
$ws.Range("G16").Value = "Chaenotheca brachypoda"
$ws.Range("H16").Value = "(Ach.) Tibell"
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("L16").ClearContents() | Out-Null
$ws.Range("P16").Value = "S om järnvägen - 8, Vg"
$ws.Range("Q16").Value = 431947.1499479365
$ws.Range("R16").Value = 6419623.056550305
$ws.Range("AJ16").Value = "tall"
$ws.Range("AK16").Value = "Pinus sylvestris"
$ws.Range("AM16").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO16").Value = "Standing dead tree/snags # Pinus sylvestris"

# Row 17
$ws.Range("A17").Value = 111483300
$ws.Range("P17").Value = "S om järnvägen - 12, Vg"
$ws.Range("Q17").Value = 431888.091041417
$ws.Range("R17").Value = 6419625.122914318
$ws.Range("AJ17").Value = "tall"
$ws.Range("AK17").Value = "Pinus sylvestris"
$ws.Range("AO17").Value = "Standing dead tree/snags # Pinus sylvestris"

# Row 18
$ws.Range("A18").Value = 111491635
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "10"
$ws.Range("K18").Value = "blomning"
$ws.Range("P18").Value = "S om järnvägen - 21, Vg"
$ws.Range("Q18").Value = 431859.6228004749
$ws.Range("R18").Value = 6419672.898494411

# Row 22
$ws.Range("A22").Value = 111661831
$ws.Range("Q22").Value = 432080.3854477856
$ws.Range("R22").Value = 6419662.773410858

# Row 23
$ws.Range("A23").Value = 111661838
$ws.Range("Q23").Value = 431799.2483237319
$ws.Range("R23").Value = 6419691.460736625
